$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 01:52"

# --- Row 4: Estados Unidos (simple value refresh) ---
$ws.Range("B4").Value = 427101
$ws.Range("C4").Value = 26766
$ws.Range("E4").Value = 390119
$ws.Range("G4").Value = 1827
$ws.Range("H4").Value = 14668

# --- Row 16: Canada (simple value refresh) ---
$ws.Range("B16").Value = 19438
$ws.Range("C16").Value = 1541
$ws.Range("E16").Value = 14463

# --- Row 17: Brasil (simple value refresh) ---
$ws.Range("B17").Value = 16188
$ws.Range("C17").Value = 2154
$ws.Range("E17").Value = 15241
$ws.Range("G17").Value = 134
$ws.Range("H17").Value = 820

# --- Rows 25-26: Australia overtakes Noruega ---
# Row 25 becomes Australia with refreshed totals
$ws.Range("A25").Value = "Australia"
$ws.Range("B25").Value = 6052
$ws.Range("C25").Value = 64
$ws.Range("D25").Value = 2813
$ws.Range("E25").Value = 3189
$ws.Range("F25").Value = 87
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 50

# Row 26 becomes Noruega, carrying its previous (unchanged) totals
$ws.Range("A26").Value = "Noruega"
$ws.Range("B26").Value = 6042
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 32
$ws.Range("E26").Value = 5909
$ws.Range("F26").Value = 78
$ws.Range("G26").Value = 12
$ws.Range("H26").Value = 101

# --- Rows 45-47: Panama overtakes Finlandia and Tailandia ---
# Row 45 becomes Panama with refreshed totals
$ws.Range("A45").Value = "Panama"
$ws.Range("B45").Value = 2528
$ws.Range("C45").Value = 428
$ws.Range("D45").Value = 16
$ws.Range("E45").Value = 2449
$ws.Range("F45").Value = 101
$ws.Range("G45").Value = 8
$ws.Range("H45").Value = 63

# Row 46 becomes Finlandia, carrying its previous (unchanged) totals
$ws.Range("A46").Value = "Finlandia"
$ws.Range("B46").Value = 2487
$ws.Range("C46").Value = 179
$ws.Range("D46").Value = 300
$ws.Range("E46").Value = 2147
$ws.Range("F46").Value = 82
$ws.Range("G46").Value = 6
$ws.Range("H46").Value = 40

# Row 47 becomes Tailandia, carrying its previous (unchanged) totals
$ws.Range("A47").Value = "Tailandia"
$ws.Range("B47").Value = 2369
$ws.Range("C47").Value = 111
$ws.Range("D47").Value = 888
$ws.Range("E47").Value = 1451
$ws.Range("F47").Value = 61
$ws.Range("G47").Value = 3
$ws.Range("H47").Value = 30

# --- Row 50: Colombia (simple value refresh) ---
$ws.Range("E50").Value = 1876
$ws.Range("G50").Value = 5
$ws.Range("H50").Value = 55

# --- Row 53: Argentina (simple value refresh) ---
$ws.Range("B53").Value = 1795
$ws.Range("C53").Value = 80
$ws.Range("E53").Value = 1372
$ws.Range("G53").Value = 5
$ws.Range("H53").Value = 65

# --- Rows 89-90: Uruguay overtakes Afganistan ---
# Row 89 becomes Uruguay with refreshed totals
$ws.Range("A89").Value = "Uruguay"
$ws.Range("B89").Value = 456
$ws.Range("C89").Value = 32
$ws.Range("D89").Value = 192
$ws.Range("E89").Value = 257
$ws.Range("F89").Value = 14
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 7

# Row 90 becomes Afganistan, carrying its previous (unchanged) totals
$ws.Range("A90").Value = "Afganistan"
$ws.Range("B90").Value = 444
$ws.Range("C90").Value = 21
$ws.Range("D90").Value = 29
$ws.Range("E90").Value = 401
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 14
